$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new contributor row at the bottom of the table (row 52)
$ws.Range("A52").Value = "Claire Sontheimer, MSW"
$ws.Range("B52").Value = "Boston University"

# Style to match the rest of the contributor table: Arial 10, wrap text,
# with left/right medium gray borders (no top/bottom) since it's the
# trailing row of the list. Apply per cell so each gets its own left+right
# edge (matching the rest of the table's per-cell border formatting).
foreach ($addr in @("A52", "B52")) {
    $cell = $ws.Range($addr)
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.WrapText = $true
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(7).Weight = -4138
    $cell.Borders.Item(7).Color = 13421772
    $cell.Borders.Item(10).LineStyle = 1
    $cell.Borders.Item(10).Weight = -4138
    $cell.Borders.Item(10).Color = 13421772
}

# Update view to match new scroll position
$ws.Application.ActiveWindow.ScrollRow = 40
